$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.459.14"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.290.64"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.21"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.15"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0960"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.338"
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.74"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.699.63"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.90"
$ws.Range("E14").Value = "  +6.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.406.99"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.284.01"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.29"
$ws.Range("E18").Value = "  +3.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.17"
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "304.54"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.41"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.97"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.36"
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.24"
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0696"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.94"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.959"
$ws.Range("E34").Value = "  +10.12%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.21"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").Value = "  +4.53%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.96"
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "126.45"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  +3.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0897"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.551"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "243.75"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.375"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0207"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.50"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.53"
$ws.Range("E51").Value = "  +1.82%  "
